# Remove leading unnamed empty column that is irrelevant to this sheet's
# test usage: shift Sheet1's data from columns B:E left to A:D, and delete
# the now-unused Sheet2 and Sheet3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 1; $r -le 4; $r++) {
  for ($c = 2; $c -le 5; $c++) {
    $srcCell = $ws.Cells.Item($r, $c)
    $dstCell = $ws.Cells.Item($r, $c - 1)
    $dstCell.Value = $srcCell.Value2
  }
}

# Clear out the now-vacated original column E (old data occupied B:E).
$ws.Range("E1:E4").Clear() | Out-Null

# Drop the two now-empty, unused worksheets.
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()
